# edit.ps1 - reproduces the commit "Add files via upload"
# Changes:
#  1) Fill in missing B/C column survey data on "Wyniki ankiet" (sheet1) rows 2-52
#  2) Bold the header rows on both "Wyniki ankiet" and "Dane ankieterów"
#  3) Adjust row heights (sheet1 row1: 45, sheet2 row1: 30) and sheet2 column F width (20)
#  4) Switch the active sheet/tab + selection from "Dane ankieterów" to "Wyniki ankiet" (cell E10)
#  5) Clear the stale selection left on "Dane ankieterów" and set its page setup (A4 portrait)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Wyniki ankiet")
$ws2 = $wb.Worksheets.Item("Dane ankieterów")

# --- 1) Fill in the missing survey results (columns B and C) on "Wyniki ankiet" ---
$ws1.Cells.Item(2, 2).Value = 4.2
$ws1.Cells.Item(2, 3).Value = 6.6
$ws1.Cells.Item(3, 2).Value = 6.5
$ws1.Cells.Item(3, 3).Value = 4.0999999999999996
$ws1.Cells.Item(4, 2).Value = 5.8
$ws1.Cells.Item(4, 3).Value = 5.4
$ws1.Cells.Item(5, 2).Value = 5.8
$ws1.Cells.Item(5, 3).Value = 5.2
$ws1.Cells.Item(6, 2).Value = 2.6
$ws1.Cells.Item(6, 3).Value = 6.1
$ws1.Cells.Item(7, 2).Value = 4.3
$ws1.Cells.Item(7, 3).Value = 6.3
$ws1.Cells.Item(8, 2).Value = 2.5
$ws1.Cells.Item(8, 3).Value = 5.8
$ws1.Cells.Item(9, 2).Value = 4.9000000000000004
$ws1.Cells.Item(9, 3).Value = 6.7
$ws1.Cells.Item(10, 2).Value = 2.2999999999999998
$ws1.Cells.Item(10, 3).Value = 5.3
$ws1.Cells.Item(11, 2).Value = 4.4000000000000004
$ws1.Cells.Item(11, 3).Value = 7.4
$ws1.Cells.Item(12, 2).Value = 5.2
$ws1.Cells.Item(12, 3).Value = 5.3
$ws1.Cells.Item(13, 2).Value = 3.6
$ws1.Cells.Item(13, 3).Value = 6.7
$ws1.Cells.Item(14, 2).Value = 6.7
$ws1.Cells.Item(14, 3).Value = 3.5
$ws1.Cells.Item(15, 2).Value = 1.8
$ws1.Cells.Item(15, 3).Value = 5.5
$ws1.Cells.Item(16, 2).Value = 5.2
$ws1.Cells.Item(16, 3).Value = 5.5
$ws1.Cells.Item(17, 2).Value = 4.5
$ws1.Cells.Item(17, 3).Value = 6.2
$ws1.Cells.Item(18, 2).Value = 3.1
$ws1.Cells.Item(18, 3).Value = 6.8
$ws1.Cells.Item(19, 2).Value = 3.7
$ws1.Cells.Item(19, 3).Value = 6.6
$ws1.Cells.Item(20, 2).Value = 3
$ws1.Cells.Item(20, 3).Value = 6.3
$ws1.Cells.Item(21, 2).Value = 4.5
$ws1.Cells.Item(21, 3).Value = 6.7
$ws1.Cells.Item(22, 2).Value = 2.9
$ws1.Cells.Item(22, 3).Value = 6.3
$ws1.Cells.Item(23, 2).Value = 4.4000000000000004
$ws1.Cells.Item(23, 3).Value = 6.2
$ws1.Cells.Item(24, 2).Value = 5.2
$ws1.Cells.Item(24, 3).Value = 4.9000000000000004
$ws1.Cells.Item(25, 2).Value = 2.7
$ws1.Cells.Item(25, 3).Value = 5.9
$ws1.Cells.Item(26, 2).Value = 3.9
$ws1.Cells.Item(26, 3).Value = 6.9
$ws1.Cells.Item(27, 2).Value = 4
$ws1.Cells.Item(27, 3).Value = 6.8
$ws1.Cells.Item(28, 2).Value = 4.5999999999999996
$ws1.Cells.Item(28, 3).Value = 6.3
$ws1.Cells.Item(29, 2).Value = 5.9
$ws1.Cells.Item(29, 3).Value = 4.5
$ws1.Cells.Item(30, 2).Value = 5.6
$ws1.Cells.Item(30, 3).Value = 4.3
$ws1.Cells.Item(31, 2).Value = 4.8
$ws1.Cells.Item(31, 3).Value = 6.2
$ws1.Cells.Item(32, 2).Value = 3.5
$ws1.Cells.Item(32, 3).Value = 6.4
$ws1.Cells.Item(33, 2).Value = 4.5999999999999996
$ws1.Cells.Item(33, 3).Value = 7.1
$ws1.Cells.Item(34, 2).Value = 4.3
$ws1.Cells.Item(34, 3).Value = 6.4
$ws1.Cells.Item(35, 2).Value = 7.1
$ws1.Cells.Item(35, 3).Value = 1.6
$ws1.Cells.Item(36, 2).Value = 5.0999999999999996
$ws1.Cells.Item(36, 3).Value = 5.9
$ws1.Cells.Item(37, 2).Value = 4.8
$ws1.Cells.Item(37, 3).Value = 5.5
$ws1.Cells.Item(38, 2).Value = 5
$ws1.Cells.Item(38, 3).Value = 4.9000000000000004
$ws1.Cells.Item(39, 2).Value = 6.3
$ws1.Cells.Item(39, 3).Value = 3.2
$ws1.Cells.Item(40, 2).Value = 6.4
$ws1.Cells.Item(40, 3).Value = 3.7
$ws1.Cells.Item(41, 2).Value = 4.7
$ws1.Cells.Item(41, 3).Value = 6
$ws1.Cells.Item(42, 2).Value = 4.0999999999999996
$ws1.Cells.Item(42, 3).Value = 6.7
$ws1.Cells.Item(43, 2).Value = 4.3
$ws1.Cells.Item(43, 3).Value = 7.1
$ws1.Cells.Item(44, 2).Value = 7.2
$ws1.Cells.Item(44, 3).Value = 2.5
$ws1.Cells.Item(45, 2).Value = 5.8
$ws1.Cells.Item(45, 3).Value = 4.3
$ws1.Cells.Item(46, 2).Value = 5.3
$ws1.Cells.Item(46, 3).Value = 4.5999999999999996
$ws1.Cells.Item(47, 2).Value = 6.5
$ws1.Cells.Item(47, 3).Value = 3.7
$ws1.Cells.Item(48, 2).Value = 7.4
$ws1.Cells.Item(48, 3).Value = 2.1
$ws1.Cells.Item(49, 2).Value = 6.1
$ws1.Cells.Item(49, 3).Value = 3.4
$ws1.Cells.Item(50, 2).Value = 7
$ws1.Cells.Item(50, 3).Value = 2
$ws1.Cells.Item(51, 2).Value = 7.2
$ws1.Cells.Item(51, 3).Value = 2.1
$ws1.Cells.Item(52, 2).Value = 4
$ws1.Cells.Item(52, 3).Value = 6.3

# --- 2) Bold the header rows ---
$ws1.Range("A1:C1").Font.Bold = $true
$ws2.Range("A1:H1").Font.Bold = $true

# --- 3) Row heights / column widths ---
$ws1.Rows.Item(1).RowHeight = 45
$ws2.Rows.Item(1).RowHeight = 30
$ws2.Columns.Item(6).ColumnWidth = 19.214285714285715

# --- 5) Page setup for "Dane ankieterów" (A4, portrait) ---
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# --- 5) Reset the stale selection on "Dane ankieterów" before moving away from it ---
$ws2.Range("A1").Select()

# --- 4) Make "Wyniki ankiet" the active sheet/tab with E10 selected ---
$ws1.Activate()
$ws1.Range("E10").Select()
